$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Update F3 value from 1400 to 1404
    $ws.Cells.Item(3, 6).Value = 1404

    # Add new row 6: copy formatting of row 5's A cell (bold/border/center style)
    # for the numbering cell A6, then set its value.
    $ws.Cells.Item(5, 1).Copy($ws.Cells.Item(6, 1))
    $ws.Cells.Item(6, 1).Value = 5

    # Column B holds a date-looking string; force text format so Excel
    # doesn't convert it into a date serial number, then clear the
    # formatting again so no stray style is left on the cell.
    $ws.Cells.Item(6, 2).NumberFormat = "@"
    $ws.Cells.Item(6, 2).Value = "2024-11-24"
    $ws.Cells.Item(6, 2).ClearFormats()

    $ws.Cells.Item(6, 3).Value = "广西·偶像梦幻祭同人ONLY"
    $ws.Cells.Item(6, 4).Value = "北湖北路48-5号(近北湖小区) 金御华尊国际大酒店"
    $ws.Cells.Item(6, 5).Value = "2024.11.24 10:00-11.24 17:00"
    $ws.Cells.Item(6, 6).Value = 2
    $ws.Cells.Item(6, 7).Value = 58
    $ws.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93736"
    $ws.Cells.Item(6, 9).Value = "//i1.hdslb.com/bfs/openplatform/202410/MsaJrW1G1728628890523.jpeg"
}
